$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.173.78'
$ws.Cells.Item(2, 5).Value = '  +0.08%  '
$ws.Cells.Item(3, 4).Value = '2.477.33'
$ws.Cells.Item(3, 5).Value = '  +0.33%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = '584.79'
$ws.Cells.Item(6, 4).Value = '174.32'
$ws.Cells.Item(6, 5).Value = '  +3.65%  '
$ws.Cells.Item(8, 4).Value = '0.514'
$ws.Cells.Item(8, 5).Value = '  +0.04%  '
$ws.Cells.Item(9, 5).Value = '  +2.50%  '
$ws.Cells.Item(10, 5).Value = '  +0.47%  '
$ws.Cells.Item(11, 5).Value = '  -0.52%  '
$ws.Cells.Item(12, 4).Value = '0.333'
$ws.Cells.Item(12, 5).Value = '  +0.59%  '
$ws.Cells.Item(13, 5).Value = '  +0.50%  '
$ws.Cells.Item(14, 4).Value = '25.45'
$ws.Cells.Item(14, 5).Value = '  -0.17%  '
$ws.Cells.Item(15, 4).Value = '67.064.35'
$ws.Cells.Item(15, 5).Value = '  +0.16%  '
$ws.Cells.Item(16, 5).Value = '  +0.29%  '
$ws.Cells.Item(17, 4).Value = '2.411.03'
$ws.Cells.Item(17, 5).Value = '  -2.45%  '
$ws.Cells.Item(18, 5).Value = '  +0.42%  '
$ws.Cells.Item(19, 4).Value = '10.95'
$ws.Cells.Item(19, 5).Value = '  -1.78%  '
$ws.Cells.Item(20, 5).Value = '  -0.96%  '
$ws.Cells.Item(21, 5).Value = '  -0.23%  '
$ws.Cells.Item(22, 5).Value = '  +0.08%  '
$ws.Cells.Item(23, 4).Value = '69.03'
$ws.Cells.Item(23, 5).Value = '  -0.20%  '
$ws.Cells.Item(24, 4).Value = '4.22'
$ws.Cells.Item(24, 5).Value = '  -0.11%  '
$ws.Cells.Item(25, 5).Value = '  +1.56%  '
$ws.Cells.Item(26, 5).Value = '  -0.66%  '
$ws.Cells.Item(27, 5).Value = '  +0.53%  '
$ws.Cells.Item(28, 5).Value = '  +0.22%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0908'
$ws.Cells.Item(29, 5).Value = '  +0.58%  '
$ws.Cells.Item(30, 4).Value = '503.07'
$ws.Cells.Item(30, 5).Value = '  -2.52%  '
$ws.Cells.Item(31, 4).Value = '7.74'
$ws.Cells.Item(31, 5).Value = '  +0.22%  '
$ws.Cells.Item(32, 5).Value = '  +0.06%  '
$ws.Cells.Item(33, 5).Value = '  -0.71%  '
$ws.Cells.Item(35, 5).Value = '  -0.08%  '
$ws.Cells.Item(36, 4).Value = '161.01'
$ws.Cells.Item(36, 5).Value = '  +0.77%  '
$ws.Cells.Item(37, 5).Value = '  +0.20%  '
$ws.Cells.Item(38, 4).Value = '18.14'
$ws.Cells.Item(38, 5).Value = '  -1.37%  '
$ws.Cells.Item(39, 4).Value = '1.33'
$ws.Cells.Item(39, 5).Value = '  -0.95%  '
$ws.Cells.Item(41, 5).Value = '  +1.64%  '
$ws.Cells.Item(42, 5).Value = '  +0.68%  '
$ws.Cells.Item(43, 5).Value = '  +1.02%  '
$ws.Cells.Item(44, 5).Value = '  +2.42%  '
$ws.Cells.Item(45, 4).Value = '142.49'
$ws.Cells.Item(45, 5).Value = '  +1.21%  '
$ws.Cells.Item(46, 5).Value = '  +0.74%  '
$ws.Cells.Item(47, 4).Value = '0.0₆0257'
$ws.Cells.Item(47, 5).Value = '  +2.33%  '
$ws.Cells.Item(48, 5).Value = '  +0.01%  '
$ws.Cells.Item(49, 4).Value = '0.0738'
$ws.Cells.Item(49, 5).Value = '  +0.29%  '
$ws.Cells.Item(50, 5).Value = '  -0.64%  '
$ws.Cells.Item(51, 5).Value = '  +0.64%  '
